# Auto-generated edit script applying market-price / profit refresh values
# per the commit diff (Sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4667.8335
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 4765.9414
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 4765.9414
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -5261.9414

$ws.Range("H67").Value = 4667.8335
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 4765.9414
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 4765.9414
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -6481.9414

$ws.Range("H74").Value = 3912.3076
$ws.Range("J74").Value = 3937.5789
$ws.Range("L74").Value = 3937.5789
$ws.Range("N74").Value = -5809.5789

$ws.Range("H76").Value = 3385.1667
$ws.Range("I76").Value = 3291.6
$ws.Range("K76").Value = 3291.6
$ws.Range("M76").Value = -2976.6

$ws.Range("H77").Value = 3912.3076
$ws.Range("J77").Value = 3937.5789
$ws.Range("L77").Value = 19687.8945
$ws.Range("N77").Value = -29047.8945

$ws.Range("H79").Value = 3385.1667
$ws.Range("I79").Value = 3291.6
$ws.Range("K79").Value = 3291.6
$ws.Range("M79").Value = -2199.6

$ws.Range("H132").Value = 502277.56
$ws.Range("I132").Value = 2090.3928
$ws.Range("J132").Value = 3503400.5
$ws.Range("K132").Value = 6271.178400000001
$ws.Range("L132").Value = 10510201.5
$ws.Range("M132").Value = -3741.178400000001
$ws.Range("N132").Value = -10515261.5

$ws.Range("H135").Value = 43940.957
$ws.Range("I135").Value = 52199.15
$ws.Range("J135").Value = 2650
$ws.Range("K135").Value = 469792.35
$ws.Range("L135").Value = 23850
$ws.Range("M135").Value = -467257.35
$ws.Range("N135").Value = -28920

$ws.Range("H137").Value = 7699037.5
$ws.Range("I137").Value = 14291800
$ws.Range("J137").Value = 7481.6665
$ws.Range("K137").Value = 42875400
$ws.Range("L137").Value = 22444.9995
$ws.Range("M137").Value = -42872850
$ws.Range("N137").Value = -27544.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 1827.5555
$ws.Range("I22").Value = 1556
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 1556
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = -1257
$ws.Range("N22").Value = -4598

$ws.Range("H32").Value = 30250.74
$ws.Range("I32").Value = 29661.656
$ws.Range("J32").Value = 31277.428
$ws.Range("K32").Value = 29661.656
$ws.Range("L32").Value = 31277.428
$ws.Range("M32").Value = -29374.656
$ws.Range("N32").Value = -31851.428

$ws.Range("H41").Value = 3472.2
$ws.Range("I41").Value = 3472.2
$ws.Range("K41").Value = 3472.2
$ws.Range("M41").Value = -3058.2

$ws.Range("H75").Value = 36448.668
$ws.Range("J75").Value = 36448.668
$ws.Range("L75").Value = 36448.668
$ws.Range("N75").Value = -38196.668

$ws.Range("H78").Value = 36448.668
$ws.Range("J78").Value = 36448.668
$ws.Range("L78").Value = 109346.004
$ws.Range("N78").Value = -118082.004

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents() | Out-Null

$ws.Range("H132").Value = 10917073
$ws.Range("I132").Value = 14316247
$ws.Range("J132").Value = 101520.73
$ws.Range("K132").Value = 42948741
$ws.Range("L132").Value = 304562.19
$ws.Range("M132").Value = -42946211
$ws.Range("N132").Value = -309622.19

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 41669210
$ws.Range("I105").Value = 62501564
$ws.Range("J105").Value = 4497.5
$ws.Range("K105").Value = 62501564
$ws.Range("L105").Value = 4497.5
$ws.Range("M105").Value = -62499817
$ws.Range("N105").Value = -7991.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3602.4
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3602.4
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3602.4
$ws.Range("M62").ClearContents() | Out-Null
$ws.Range("N62").Value = -4850.4

$ws.Range("H65").Value = 3602.4
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3602.4
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 18012
$ws.Range("M65").ClearContents() | Out-Null
$ws.Range("N65").Value = -24252

$ws.Range("H86").Value = 1994
$ws.Range("I86").Value = 1310.2
$ws.Range("J86").Value = 3361.6
$ws.Range("K86").Value = 1310.2
$ws.Range("L86").Value = 3361.6
$ws.Range("M86").Value = -187.2
$ws.Range("N86").Value = -5607.6

$ws.Range("H89").Value = 1994
$ws.Range("I89").Value = 1310.2
$ws.Range("J89").Value = 3361.6
$ws.Range("K89").Value = 6551
$ws.Range("L89").Value = 16808
$ws.Range("M89").Value = -935
$ws.Range("N89").Value = -28040

$ws.Range("H99").Value = 1633.3334
$ws.Range("I99").Value = 1450
$ws.Range("K99").Value = 1450
$ws.Range("M99").Value = 48

$ws.Range("H126").Value = 1633.3334
$ws.Range("I126").Value = 1450
$ws.Range("K126").Value = 4350
$ws.Range("M126").Value = -1880

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 797.7826
$ws.Range("J34").Value = 958.7646999999999
$ws.Range("L34").Value = 2876.2941
$ws.Range("N34").Value = -3044.2941

$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents() | Out-Null

$ws.Range("H55").Value = 2250
$ws.Range("J55").Value = 2600
$ws.Range("L55").Value = 7800
$ws.Range("N55").Value = -8154

$ws.Range("H87").Value = 31075
$ws.Range("I87").Value = 24533.334
$ws.Range("K87").Value = 73600.00199999999
$ws.Range("M87").Value = -72352.00199999999

$ws.Range("H90").Value = 31075
$ws.Range("I90").Value = 24533.334
$ws.Range("K90").Value = 220800.006
$ws.Range("M90").Value = -214560.006

$ws.Range("H139").Value = 2154.1667
$ws.Range("I139").Value = 1810.2941
$ws.Range("J139").Value = 8000
$ws.Range("K139").Value = 5430.8823
$ws.Range("L139").Value = 24000
$ws.Range("M139").Value = -290.8823000000002
$ws.Range("N139").Value = -34280

$ws.Range("H141").Value = 7461.4287
$ws.Range("I141").Value = 7461.4287
$ws.Range("K141").Value = 22384.2861
$ws.Range("M141").Value = -17204.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 50608.863
$ws.Range("I70").Value = 76363.92999999999
$ws.Range("J70").Value = 5537.5
$ws.Range("K70").Value = 76363.92999999999
$ws.Range("L70").Value = 5537.5
$ws.Range("M70").Value = -76093.92999999999
$ws.Range("N70").Value = -6077.5

$ws.Range("H73").Value = 50608.863
$ws.Range("I73").Value = 76363.92999999999
$ws.Range("J73").Value = 5537.5
$ws.Range("K73").Value = 76363.92999999999
$ws.Range("L73").Value = 5537.5
$ws.Range("M73").Value = -75427.92999999999
$ws.Range("N73").Value = -7409.5

$ws.Range("H80").Value = 3809.658
$ws.Range("J80").Value = 3777.7878
$ws.Range("L80").Value = 3777.7878
$ws.Range("N80").Value = -5773.7878

$ws.Range("H83").Value = 3809.658
$ws.Range("J83").Value = 3777.7878
$ws.Range("L83").Value = 18888.939
$ws.Range("N83").Value = -28872.939

$ws.Range("H123").Value = 22896.941
$ws.Range("J123").Value = 22896.941
$ws.Range("L123").Value = 22896.941
$ws.Range("N123").Value = -27796.941

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 50001980
$ws.Range("I7").Value = 62502100
$ws.Range("J7").Value = 1500
$ws.Range("K7").Value = 62502100
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = -62501988
$ws.Range("N7").Value = -1724

$ws.Range("H22").Value = 941.5
$ws.Range("I22").Value = 878.5833
$ws.Range("J22").Value = 1067.3334
$ws.Range("K22").Value = 878.5833
$ws.Range("L22").Value = 1067.3334
$ws.Range("M22").Value = -583.5833
$ws.Range("N22").Value = -1657.3334

$ws.Range("H27").Value = 941.5
$ws.Range("I27").Value = 878.5833
$ws.Range("J27").Value = 1067.3334
$ws.Range("K27").Value = 878.5833
$ws.Range("L27").Value = 1067.3334
$ws.Range("M27").Value = -771.5833
$ws.Range("N27").Value = -1281.3334

$ws.Range("H40").Value = 2855.889
$ws.Range("I40").Value = 2783
$ws.Range("K40").Value = 2783
$ws.Range("M40").Value = -2647

$ws.Range("H46").Value = 758565.5
$ws.Range("I46").Value = 1318017.9
$ws.Range("J46").Value = 1659.2941
$ws.Range("K46").Value = 1318017.9
$ws.Range("L46").Value = 1659.2941
$ws.Range("M46").Value = -1317829.9
$ws.Range("N46").Value = -2035.2941

$ws.Range("H126").Value = 50001980
$ws.Range("I126").Value = 62502100
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 187506300
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -187503830
$ws.Range("N126").Value = -9440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2709.9
$ws.Range("I81").Value = 1800.5
$ws.Range("K81").Value = 3601
$ws.Range("M81").Value = -2540

$ws.Range("H84").Value = 2709.9
$ws.Range("I84").Value = 1800.5
$ws.Range("K84").Value = 18005
$ws.Range("M84").Value = -12701

$ws.Range("H122").Value = 2056.25
$ws.Range("I122").Value = 1620
$ws.Range("J122").Value = 2254.5454
$ws.Range("K122").Value = 4860
$ws.Range("L122").Value = 6763.6362
$ws.Range("M122").Value = -2410
$ws.Range("N122").Value = -11663.6362

$ws.Range("H132").Value = 38444.203
$ws.Range("I132").Value = 30404.824
$ws.Range("J132").Value = 52111.15
$ws.Range("K132").Value = 91214.47200000001
$ws.Range("L132").Value = 156333.45
$ws.Range("M132").Value = -88684.47200000001
$ws.Range("N132").Value = -161393.45
